$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 411
$ws.Range("J2").Value = 1688
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 463
$ws.Range("M2").Value = 27
$ws.Range("N2").Value = 273
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 6
$ws.Range("R2").Value = 24
$ws.Range("S2").Value = 182
$ws.Range("T2").Value = 273
$ws.Range("U2").Value = 24
$ws.Range("V2").Value = 2693
$ws.Range("X2").Value = 2593
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 35
$ws.Range("AA2").Value = 19
